# Daily attendance processing - 2026-01-29 19:46:42
# Normalises the "Recorded By" (column G) value on the Session Analysis
# Results sheet: the comma-separated list of recorder names/emails is
# re-ordered (alphabetising "System" vs "system" and ordering e-mails
# before the bare "System" sentinel etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact text substitutions observed for the "Recorded By" column.
$replacements = @{
    "dnasr281@gmail.com, System"              = "System, dnasr281@gmail.com"
    "backup@backdoor.com, system, System"     = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, admin@admin.com"     = "admin@admin.com, dnasr281@gmail.com"
}

# Find the last used row on the sheet and walk column G (7), re-writing
# any cell whose text matches one of the known "before" values.
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    if ($null -ne $current -and $replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
